$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = 0.86740239339518777
$ws.Range("BK2").Value = 0.92701522327619634
$ws.Range("A3").Value = 0.93722787509806649
$ws.Range("E3").Value = 0.67184613495673418
$ws.Range("Z3").Value = 0.81991840624012657
$ws.Range("C4").Value = 0.85543738486878707
$ws.Range("E4").Value = 0.95622892660015402
$ws.Range("F4").Value = 0.98427473306925795
$ws.Range("N5").Value = 0.91326143795681158
$ws.Range("G6").Value = 0.63503441034154129
$ws.Range("B7").Value = 0.92935077059209159
$ws.Range("E7").Value = 0.98616252439679308
$ws.Range("F8").Value = 0.88120263445916347
$ws.Range("G8").Value = 0.72441312751801568
$ws.Range("I8").Value = 0.52334954534462907
$ws.Range("G9").Value = 0.90176664548076046
$ws.Range("H10").Value = 0.85056422623884265
$ws.Range("I10").Value = 0.64305046333478288
$ws.Range("K10").Value = 0.99613024542701223
$ws.Range("L10").Value = 0.63743685871053168
$ws.Range("I11").Value = 0.79789926530792576
$ws.Range("T11").Value = 0.81245075289144841
$ws.Range("K12").Value = 0.90299395240326219
$ws.Range("K13").Value = 0.64496635352592213
$ws.Range("L13").Value = 0.86525886916928985
$ws.Range("W13").Value = 0.99771497966381628
$ws.Range("L14").Value = 0.97641571714225051
$ws.Range("M14").Value = 0.98217819250784122
$ws.Range("P14").Value = 0.7989771311656203
$ws.Range("AT14").Value = 0.82940656137408242
$ws.Range("M15").Value = 0.78467279854210881
$ws.Range("N15").Value = 0.86146825822608775
$ws.Range("P15").Value = 0.8410116085240309
$ws.Range("Q15").Value = 0.91578256637838673
$ws.Range("R16").Value = 0.87932664416782247
$ws.Range("BC16").Value = 0.92557720079560313
$ws.Range("S17").Value = 0.94126944770436249
$ws.Range("Q18").Value = 0.94541371900446713
$ws.Range("AP18").Value = 0.8725603943041218
$ws.Range("R19").Value = 0.79883161319474971
$ws.Range("V19").Value = 0.96712708720464002
$ws.Range("Y20").Value = 0.87648938525850062
$ws.Range("T21").Value = 0.6232842992113159
$ws.Range("U22").Value = 0.94955408323488699
$ws.Range("W22").Value = 0.77981590113819843
$ws.Range("X22").Value = 0.61220889509913823
$ws.Range("N23").Value = 0.61189477413833759
$ws.Range("U23").Value = 0.60445379004906163
$ws.Range("P24").Value = 0.86820098985602234
$ws.Range("Y24").Value = 0.9980501910018611
$ws.Range("AB27").Value = 0.78024172359698529
$ws.Range("AC27").Value = 0.91243446708133003
$ws.Range("S29").Value = 0.65502801248730369
$ws.Range("AB30").Value = 0.99261445222450839
$ws.Range("AC30").Value = 0.81684660662940323
$ws.Range("AB31").Value = 0.71615630231735783
$ws.Range("AC31").Value = 0.87933856514370345
$ws.Range("AD31").Value = 0.83402955824088554
$ws.Range("Z32").Value = 0.78456532771741072
$ws.Range("AD32").Value = 0.79487719718357175
$ws.Range("AE32").Value = 0.74774099980268416
$ws.Range("AG32").Value = 0.87296416625057072
$ws.Range("AJ32").Value = 0.88726388685179902
$ws.Range("AH33").Value = 0.84101106119879132
$ws.Range("BF33").Value = 0.94423344620319161
$ws.Range("AF34").Value = 0.98316508805926861
$ws.Range("AI34").Value = 0.74125751574776677
$ws.Range("AQ34").Value = 0.80751400024692077
$ws.Range("BN34").Value = 0.64617234765734044
$ws.Range("AG35").Value = 0.55829787128022634
$ws.Range("BN35").Value = 0.61878356578817817
$ws.Range("AI36").Value = 0.66771600722470748
$ws.Range("AI37").Value = 0.92877015593254453
$ws.Range("AL37").Value = 0.5850766578661758
$ws.Range("AM37").Value = 0.73807309235473939
$ws.Range("BK37").Value = 0.93793235923304541
$ws.Range("AJ38").Value = 0.81269281761759382
$ws.Range("AN38").Value = 0.67593221378135981
$ws.Range("BL38").Value = 0.75683761243475023
$ws.Range("AL39").Value = 0.92351232955530094
$ws.Range("BF39").Value = 0.86091186359083594
$ws.Range("AP40").Value = 0.74158838543627881
$ws.Range("AM41").Value = 0.66486725922874901
$ws.Range("AN41").Value = 0.78562686471828425
$ws.Range("AQ41").Value = 0.99866059827701359
$ws.Range("AD42").Value = 0.9972574036435462
$ws.Range("AO42").Value = 0.98465188803309067
$ws.Range("AR42").Value = 0.59709194405555044
$ws.Range("Y43").Value = 0.65552708742175503
$ws.Range("BK43").Value = 0.66337187239363882
$ws.Range("AP45").Value = 0.96709262795030537
$ws.Range("AQ45").Value = 0.83521692739218012
$ws.Range("AR45").Value = 0.62385884510891398
$ws.Range("AX45").Value = 0.80165398160104773
$ws.Range("AR46").Value = 0.90133384371731573
$ws.Range("N47").Value = 0.8127254425442173
$ws.Range("AX47").Value = 0.94868466269389384
$ws.Range("AT48").Value = 0.69081036207489932
$ws.Range("AU49").Value = 0.6921813395553218
$ws.Range("AV49").Value = 0.75214697136221054
$ws.Range("AX49").Value = 0.71342519592618847
$ws.Range("AV50").Value = 0.93556732386248975
$ws.Range("AW51").Value = 0.92804154008044648
$ws.Range("AZ51").Value = 0.96069922649680839
$ws.Range("BA51").Value = 0.93625236409024604
$ws.Range("AZ53").Value = 0.68743585822039166
$ws.Range("AZ54").Value = 0.85639515920703579
$ws.Range("BA54").Value = 0.80074732912500568
$ws.Range("BC54").Value = 0.66057553619243559
$ws.Range("BD54").Value = 0.73224365861320995
$ws.Range("BA55").Value = 0.98469557628247806
$ws.Range("BC56").Value = 0.71423783790969242
$ws.Range("BE56").Value = 0.8980615381748569
$ws.Range("BF56").Value = 0.93444283044573617
$ws.Range("BC57").Value = 0.7634863441585571
$ws.Range("BG57").Value = 0.72846394610202392
$ws.Range("AX58").Value = 0.72977245144525493
$ws.Range("BE58").Value = 0.80688620642606934
$ws.Range("BH58").Value = 0.90620950368919839
$ws.Range("BH59").Value = 0.92651498247517594
$ws.Range("BI60").Value = 0.97787916735679392
$ws.Range("BJ60").Value = 0.73893096427070981
$ws.Range("BG61").Value = 0.97876389468691305
$ws.Range("BI62").Value = 0.86856815579248936
$ws.Range("BI63").Value = 0.70774189897667927
$ws.Range("BJ64").Value = 0.74248375411772038
$ws.Range("BK64").Value = 0.71991660003647628
$ws.Range("BM64").Value = 0.77520769345322638
$ws.Range("BJ65").Value = 0.79712068863079388
$ws.Range("BM66").Value = 0.95679326541115284
$ws.Range("AR67").Value = 0.94538368002903839
$ws.Range("BM67").Value = 0.96581974988098274
$ws.Range("BP67").Value = 0.98778023373866053
$ws.Range("A68").Value = 0.8728544145123851
$ws.Range("AB68").Value = 0.77178991545711673
$ws.Range("AN68").Value = 0.98814870016466205
